# PGM_04.pptx — lecture 05, change final date of seminar project
#
# 1) Delete the trailing "Kahoot time!" slide (last slide, slide 10).
# 2) Refresh the cached "datetimeFigureOut" footer date (slide master +
#    every slide layout) from 10/12/2020 to 11/8/2021.

$p = $ppt.ActivePresentation

# --- 1) Remove the "Kahoot time!" slide -----------------------------------
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    $isKahoot = $false
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "Kahoot time!") {
                $isKahoot = $true
            }
        }
    }
    if ($isKahoot) {
        $slide.Delete()
    }
}

# --- 2) Update the cached date placeholder text ----------------------------
$newDate = "11/8/2021"

function Update-DatePlaceholder($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.Type -eq 14) {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                if ($shp.HasTextFrame) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}
